# Update "想去人数" (interested-people count) figures that changed between
# the previous and newly generated data snapshot.
#
# Sheet "展览" (Exhibitions):
#   F3: 205 -> 206   (南宁·熊喵M动漫嘉年华·万圣派对)
#   F4: 816 -> 818   (南宁·万圣漫控嘉年华10)
#
# Sheet "全部类型" (All types) contains the same rows shifted down by one
# (it has an extra leading data row), so the matching cells are F4/F5:
#   F4: 205 -> 206
#   F5: 816 -> 818

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 206
$wsExhibition.Range("F4").Value = 818

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 206
$wsAllTypes.Range("F5").Value = 818
